$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.249.07'
$ws.Range('E2').Value = '  +0.28%  '

$ws.Range('D3').Value = '1.687.17'
$ws.Range('E3').Value = '  +0.93%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.72'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5235'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.26%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2690'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06441'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.98%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.09'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.10%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07463'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.02%  '

$ws.Range('D12').Value = '1.693.82'
$ws.Range('E12').Value = '  +1.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.547'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.29%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5860'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.89%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008559'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.34%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.77'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.27%  '

$ws.Range('D17').Value = '26.319.03'
$ws.Range('E17').Value = '  +0.35%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.977'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.37%  '

$ws.Range('E19').Value = '  -0.12%  '

$ws.Range('E20').Value = '  +0.46%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.89'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.29%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.244'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.05%  '

$ws.Range('E23').Value = '  -0.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.30'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.54%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1247'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +7.00%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.655'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.38%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.06874'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +19.73%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.87'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.40%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.343'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.58%  '

$ws.Range('E30').Value = '  -0.15%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.596'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.72%  '

$ws.Range('E32').Value = '  +1.79%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.660'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.39%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.027'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.27%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6204'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.81%  '

$ws.Range('E36').Value = '  +0.25%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.707'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.81%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.298'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +6.64%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01619'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.31%  '

$ws.Range('D40').Value = '1.100.95'
$ws.Range('E40').Value = '  +0.94%  '

$ws.Range('E41').Value = '  +1.70%  '

$ws.Range('E42').Value = '  +0.97%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.97'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.15%  '

$ws.Range('D44').Value = '1.836.82'
$ws.Range('E44').Value = '  +0.95%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000111'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.94%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.94'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.65%  '

$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.008'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.55%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.160'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.30%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05256'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.06%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4289'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.60%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.008'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.34%  '
